$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.149.76"
$ws.Range("E2").Value = "  -1.26%  "

$ws.Range("D3").Value = "2.288.97"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.36%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.12"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  -0.24%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.584"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.97%  "

$ws.Range("D9").Value = "2.286.94"
$ws.Range("E9").Value = "  -0.83%  "

$ws.Range("E10").Value = "  -3.11%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.32%  "

$ws.Range("E12").Value = "  +0.14%  "

$ws.Range("E13").Value = "  -2.44%  "

$ws.Range("E14").Value = "  -1.86%  "

$ws.Range("D15").Value = "2.696.10"
$ws.Range("E15").Value = "  -1.38%  "

$ws.Range("D16").Value = "58.040.41"
$ws.Range("E16").Value = "  -1.50%  "

$ws.Range("E17").Value = "  -1.68%  "

$ws.Range("D18").Value = "2.282.88"
$ws.Range("E18").Value = "  -1.30%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.52"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.88%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "311.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.86%  "

$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "62.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.44%  "

$ws.Range("E25").Value = "  -0.88%  "

$ws.Range("E26").Value = "  -0.31%  "

$ws.Range("E27").Value = "  -3.63%  "

$ws.Range("E28").Value = "  -5.90%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.96"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.16%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.71"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.12%  "

$ws.Range("D31").Value = "0.0₃0718"
$ws.Range("E31").Value = "  -2.51%  "

$ws.Range("E32").Value = "  -2.50%  "

$ws.Range("E33").Value = "  -4.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.376"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.94%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.72"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("E37").Value = "  +0.08%  "

$ws.Range("E38").Value = "  -4.61%  "

$ws.Range("E39").Value = "  -3.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.28%  "

$ws.Range("E41").Value = "  -4.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.57"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.05%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "287.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.56%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.40"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.23%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0948"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.22%  "

$ws.Range("E46").Value = "  -1.01%  "

$ws.Range("E47").Value = "  -0.83%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.02"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.28%  "

$ws.Range("E49").Value = "  -2.25%  "

$ws.Range("E50").Value = "  -0.92%  "

$ws.Range("E51").Value = "  -0.73%  "
